# Generate Report for Handoff
#
# Updates the localization-status report after a handoff/handback run:
#  - Overview sheet: "Latest HO Xliff Generate Date" (column G) refreshed
#    for the six files that were just handed off.
#  - de-de sheet: "Latest Handoff Datetime" (column H) refreshed for the
#    same six files (it shares the same generate timestamp as Overview).
#  - zh-cn sheet: "Latest Handoff Datetime" (column H) refreshed with its
#    own (slightly earlier) generate timestamp for the same six files.
#  - zh-cn and de-de sheets: "Priority" (column E) set to "ht" for the
#    same six files, now that the handoff type has been determined.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$overviewGenerateDate = "2016-09-04 04:24:13"
$zhCnHandoffDate      = "2016-09-04 04:24:08"
$deDeHandoffDate      = "2016-09-04 04:24:13"
$priority             = "ht"

# Rows 7-12 correspond to the six files whose handoff report is being
# regenerated: 5313fc18, 72f7cac2, 7fbfb5c3, 90545bff, 9a5a8730, b0c3bfe8.
foreach ($r in 7..12) {
    $wsOverview.Range("G$r").Value = $overviewGenerateDate

    $wsZhCn.Range("H$r").Value = $zhCnHandoffDate
    $wsZhCn.Range("E$r").Value = $priority

    $wsDeDe.Range("H$r").Value = $deDeHandoffDate
    $wsDeDe.Range("E$r").Value = $priority
}
